$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated price (Price) and volume (Volume(1h)) figures refreshed by the crypto-price scraper.
# Cells are stored as text in the workbook (inlineStr), so force a Text number format
# before assigning the value to stop Excel from re-interpreting numeric-looking strings
# or percentages as real numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "320.02"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.97%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "49.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.13%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.337"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.77%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08025"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.33%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.603"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.00%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.342"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "26.95%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.651"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.89%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.81%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1966"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.99%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09610"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "4.64%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04530"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.90%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.37%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001309"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.39%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.04217"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.04%"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.19%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.60%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.477"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "5.70%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.88%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.160"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.14%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1393"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.91%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.3015"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.47%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001294"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.45%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004304"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-4.06%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.63%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003546"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02684"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "3.99%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05918"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "10.68%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "92.84%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008040"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.06%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1468"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.32%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007528"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.02%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007923"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.18%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3216"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.54%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007027"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.00%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05537"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-15.75%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004007"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.83%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
